$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataBF = New-Object 'object[,]' 24,5
$dataBF[0,0] = 1.02
$dataBF[0,1] = 1.01221018469611
$dataBF[0,2] = 1.045569785128762
$dataBF[0,3] = 1.014288546000627
$dataBF[0,4] = 1.048685544409157
$dataBF[1,0] = 1.02
$dataBF[1,1] = 1.013116714357632
$dataBF[1,2] = 1.046147893768083
$dataBF[1,3] = 1.015055512287565
$dataBF[1,4] = 1.049446784146355
$dataBF[2,0] = 1.02
$dataBF[2,1] = 1.013703795212341
$dataBF[2,2] = 1.046520266251091
$dataBF[2,3] = 1.015552622526932
$dataBF[2,4] = 1.049937915616048
$dataBF[3,0] = 1.02
$dataBF[3,1] = 1.013950721462018
$dataBF[3,2] = 1.046676401562198
$dataBF[3,3] = 1.015761805262442
$dataBF[3,4] = 1.050144039545013
$dataBF[4,0] = 1.02
$dataBF[4,1] = 1.013992188319462
$dataBF[4,2] = 1.046702593247039
$dataBF[4,3] = 1.015796939506255
$dataBF[4,4] = 1.050178628190763
$dataBF[5,0] = 1.02
$dataBF[5,1] = 1.013707094192926
$dataBF[5,2] = 1.046522354153342
$dataBF[5,3] = 1.015555416862429
$dataBF[5,4] = 1.049940671223869
$dataBF[6,0] = 1.02
$dataBF[6,1] = 1.01251644732684
$dataBF[6,2] = 1.045765510659034
$dataBF[6,3] = 1.014547572703912
$dataBF[6,4] = 1.048943105850657
$dataBF[7,0] = 1.02
$dataBF[7,1] = 1.010422222433714
$dataBF[7,2] = 1.044418926061806
$dataBF[7,3] = 1.012778060446377
$dataBF[7,4] = 1.047174342393254
$dataBF[8,0] = 1.02
$dataBF[8,1] = 1.009028727386166
$dataBF[8,2] = 1.043512654874704
$dataBF[8,3] = 1.011602804547341
$dataBF[8,4] = 1.045987973050702
$dataBF[9,0] = 1.02
$dataBF[9,1] = 1.008425972335102
$dataBF[9,2] = 1.043118231916976
$dataBF[9,3] = 1.011094972761635
$dataBF[9,4] = 1.045472587167383
$dataBF[10,0] = 1.02
$dataBF[10,1] = 1.008202179127635
$dataBF[10,2] = 1.042971427313497
$dataBF[10,3] = 1.010906502454154
$dataBF[10,4] = 1.045280899943204
$dataBF[11,0] = 1.02
$dataBF[11,1] = 1.008250179130058
$dataBF[11,2] = 1.043002930863002
$dataBF[11,3] = 1.010946922663809
$dataBF[11,4] = 1.045322028775497
$dataBF[12,0] = 1.02
$dataBF[12,1] = 1.008407471536686
$dataBF[12,2] = 1.043106103075512
$dataBF[12,3] = 1.011079390453574
$dataBF[12,4] = 1.045456747326622
$dataBF[13,0] = 1.02
$dataBF[13,1] = 1.008504397446361
$dataBF[13,2] = 1.04316963141567
$dataBF[13,3] = 1.01116102961379
$dataBF[13,4] = 1.045539718823229
$dataBF[14,0] = 1.02
$dataBF[14,1] = 1.009068743777289
$dataBF[14,2] = 1.043538789458782
$dataBF[14,3] = 1.011636530169779
$dataBF[14,4] = 1.046022142326786
$dataBF[15,0] = 1.02
$dataBF[15,1] = 1.009422914759754
$dataBF[15,2] = 1.043769818346784
$dataBF[15,3] = 1.011935084440361
$dataBF[15,4] = 1.046324305774273
$dataBF[16,0] = 1.02
$dataBF[16,1] = 1.009629558138939
$dataBF[16,2] = 1.043904380311411
$dataBF[16,3] = 1.012109328440131
$dataBF[16,4] = 1.046500390329577
$dataBF[17,0] = 1.02
$dataBF[17,1] = 1.009700028579078
$dataBF[17,2] = 1.043950229608683
$dataBF[17,3] = 1.012168758454653
$dataBF[17,4] = 1.046560403013136
$dataBF[18,0] = 1.02
$dataBF[18,1] = 1.009384909202592
$dataBF[18,2] = 1.043745051107941
$dataBF[18,3] = 1.011903041795676
$dataBF[18,4] = 1.046291903262571
$dataBF[19,0] = 1.02
$dataBF[19,1] = 1.008361150160978
$dataBF[19,2] = 1.043075729647737
$dataBF[19,3] = 1.011040377533879
$dataBF[19,4] = 1.045417082951586
$dataBF[20,0] = 1.02
$dataBF[20,1] = 1.007718033132291
$dataBF[20,2] = 1.042653175631372
$dataBF[20,3] = 1.010498918680515
$dataBF[20,4] = 1.044865603707019
$dataBF[21,0] = 1.02
$dataBF[21,1] = 1.008058908109133
$dataBF[21,2] = 1.042877342247162
$dataBF[21,3] = 1.01078586737943
$dataBF[21,4] = 1.045158089444237
$dataBF[22,0] = 1.02
$dataBF[22,1] = 1.009402082088177
$dataBF[22,2] = 1.043756242954464
$dataBF[22,3] = 1.011917520171477
$dataBF[22,4] = 1.046306545064424
$dataBF[23,0] = 1.02
$dataBF[23,1] = 1.010963166826098
$dataBF[23,2] = 1.04476856761513
$dataBF[23,3] = 1.013234749512647
$dataBF[23,4] = 1.047632889101172

$dataIN = New-Object 'object[,]' 24,6
$dataIN[0,0] = 1.039436976155305
$dataIN[0,1] = 1.017454211236662
$dataIN[0,2] = 1.048337655675921
$dataIN[0,3] = 1.017147146308949
$dataIN[0,4] = 1.05144469945639
$dataIN[0,5] = 1.00994939178592
$dataIN[1,0] = 1.039560038532723
$dataIN[1,1] = 1.017994566395177
$dataIN[1,2] = 1.048727839110871
$dataIN[1,3] = 1.017719132122746
$dataIN[1,4] = 1.052018165815712
$dataIN[1,5] = 1.010129820119729
$dataIN[2,0] = 1.039637156974232
$dataIN[2,1] = 1.018344121375491
$dataIN[2,2] = 1.048977914197344
$dataIN[2,3] = 1.018089420633729
$dataIN[2,4] = 1.0523871363357
$dataIN[2,5] = 1.01024647998461
$dataIN[3,0] = 1.039668975046893
$dataIN[3,1] = 1.018491051407282
$dataIN[3,2] = 1.049082469114757
$dataIN[3,3] = 1.018245130742215
$dataIN[3,4] = 1.052541746206001
$dataIN[3,5] = 1.010295501992123
$dataIN[4,0] = 1.039674282069646
$dataIN[4,1] = 1.018515720231938
$dataIN[4,2] = 1.049099990467059
$dataIN[4,3] = 1.018271277513722
$dataIN[4,4] = 1.052567676178244
$dataIN[4,5] = 1.010303731712491
$dataIN[5,0] = 1.039637584498116
$dataIN[5,1] = 1.018346084753088
$dataIN[5,2] = 1.048979313533907
$dataIN[5,3] = 1.01809150108093
$dataIN[5,4] = 1.05238920422917
$dataIN[5,5] = 1.010247135105255
$dataIN[6,0] = 1.039479084807838
$dataIN[6,1] = 1.017636844955866
$dataIN[6,2] = 1.048470015577145
$dataIN[6,3] = 1.017340414556784
$dataIN[6,4] = 1.05163893856484
$dataIN[6,5] = 1.01001038663488
$dataIN[7,0] = 1.039180614189156
$dataIN[7,1] = 1.016386426094656
$dataIN[7,2] = 1.047554282607424
$dataIN[7,3] = 1.016018299797517
$dataIN[7,4] = 1.050300893789225
$dataIN[7,5] = 1.009592539567516
$dataIN[8,0] = 1.038968819703922
$dataIN[8,1] = 1.015552442958052
$dataIN[8,2] = 1.046931624623644
$dataIN[8,3] = 1.015137897706406
$dataIN[8,4] = 1.049398261267985
$dataIN[8,5] = 1.009313551294094
$dataIN[9,0] = 1.038874086912159
$dataIN[9,1] = 1.015191244599724
$dataIN[9,2] = 1.046659146050773
$dataIN[9,3] = 1.014756927325849
$dataIN[9,4] = 1.049004925866656
$dataIN[9,5] = 1.009192650858092
$dataIN[10,0] = 1.038838445778189
$dataIN[10,1] = 1.015057068759728
$dataIN[10,2] = 1.04655750715481
$dataIN[10,3] = 1.014615456594588
$dataIN[10,4] = 1.048858451842831
$dataIN[10,5] = 1.009147728924999
$dataIN[11,0] = 1.038846111409308
$dataIN[11,1] = 1.015085850399616
$dataIN[11,2] = 1.046579328389236
$dataIN[11,3] = 1.01464580078028
$dataIN[11,4] = 1.048889887803905
$dataIN[11,5] = 1.009157365465126
$dataIN[12,0] = 1.038871150042181
$dataIN[12,1] = 1.015180153794869
$dataIN[12,2] = 1.046650753272098
$dataIN[12,3] = 1.014745232516433
$dataIN[12,4] = 1.048992825848484
$dataIN[12,5] = 1.009188937883015
$dataIN[13,0] = 1.038886517158765
$dataIN[13,1] = 1.01523825584438
$dataIN[13,2] = 1.046694703815339
$dataIN[13,3] = 1.014806500843742
$dataIN[13,4] = 1.049056200183476
$dataIN[13,5] = 1.009208388831403
$dataIN[14,0] = 1.03897504319183
$dataIN[14,1] = 1.015576412946992
$dataIN[14,2] = 1.046949647968287
$dataIN[14,3] = 1.015163186810779
$dataIN[14,4] = 1.049424313386786
$dataIN[14,5] = 1.009321573054909
$dataIN[15,0] = 1.039029764389669
$dataIN[15,1] = 1.015788509784907
$dataIN[15,2] = 1.047108802462652
$dataIN[15,3] = 1.015386994108909
$dataIN[15,4] = 1.049654556171133
$dataIN[15,5] = 1.009392544926084
$dataIN[16,0] = 1.039061390448683
$dataIN[16,1] = 1.015912214670602
$dataIN[16,2] = 1.04720135816268
$dataIN[16,3] = 1.015517561119598
$dataIN[16,4] = 1.049788612533481
$dataIN[16,5] = 1.009433932263839
$dataIN[17,0] = 1.039072124570461
$dataIN[17,1] = 1.015954393516608
$dataIN[17,2] = 1.047232870317112
$dataIN[17,3] = 1.015562085132161
$dataIN[17,4] = 1.049834281449764
$dataIN[17,5] = 1.009448042677806
$dataIN[18,0] = 1.039023923505642
$dataIN[18,1] = 1.015765754581096
$dataIN[18,2] = 1.047091755268902
$dataIN[18,3] = 1.015362979213183
$dataIN[18,4] = 1.0496298781441
$dataIN[18,5] = 1.00938493128063
$dataIN[19,0] = 1.038863789288508
$dataIN[19,1] = 1.015152384081994
$dataIN[19,2] = 1.04662973222946
$dataIN[19,3] = 1.014715951272867
$dataIN[19,4] = 1.048962523403725
$dataIN[19,5] = 1.009179640979702
$dataIN[20,0] = 1.038760485313843
$dataIN[20,1] = 1.014766672389587
$dataIN[20,2] = 1.046336763412884
$dataIN[20,3] = 1.014309363193648
$dataIN[20,4] = 1.048540781312764
$dataIN[20,5] = 1.009050485273617
$dataIN[21,0] = 1.038815496765464
$dataIN[21,1] = 1.014971150916646
$dataIN[21,2] = 1.046492305820005
$dataIN[21,3] = 1.014524881608569
$dataIN[21,4] = 1.048764557784398
$dataIN[21,5] = 1.009118960755002
$dataIN[22,0] = 1.039026563651688
$dataIN[22,1] = 1.015776036701817
$dataIN[22,2] = 1.047099459016215
$dataIN[22,3] = 1.015373830436801
$dataIN[22,4] = 1.049641029825086
$dataIN[22,5] = 1.009388371588224
$dataIN[23,0] = 1.03926004015342
$dataIN[23,1] = 1.016709759390862
$dataIN[23,2] = 1.047793175446224
$dataIN[23,3] = 1.016359925592609
$dataIN[23,4] = 1.050648689496398
$dataIN[23,5] = 1.009700639472655

$ws.Range("B2:F25").Value = $dataBF
$ws.Range("I2:N25").Value = $dataIN
